$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "edit1"
$ws.Range("B13").Value = "riya-morankar"
$ws.Range("C13").Value = "Squashed"
$ws.Range("D13").Value = "N/A"
$ws.Range("F13").Value = "846895ac8c5fcf9bec1e93cf92bcdf081ff57046"

# E13 holds a date-shaped string ("2025-06-18") that must stay plain text,
# matching the other Date-column cells in this sheet (not an actual date
# serial value). Force text formatting before typing it so Excel doesn't
# auto-convert it to a date, then pull the plain/default formatting back in
# from a neighboring already-text cell so no extra date-number-format style
# lingers on the cell.
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2025-06-18"
$ws.Range("D2").Copy()
$ws.Range("E13").PasteSpecial(-4122)
